$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 355.375  # ALC!H2: 391.33334 -> 355.375
$ws.Cells.Item(2, 9).Value = 396.33334  # ALC!I2: 391.33334 -> 396.33334
$ws.Cells.Item(2, 10).Value = 232.5  # ALC!J2: 0 -> 232.5
$ws.Cells.Item(2, 11).Value = 396.33334  # ALC!K2: 391.33334 -> 396.33334
$ws.Cells.Item(2, 12).Value = 232.5  # ALC!L2: 0 -> 232.5
$ws.Cells.Item(2, 13).Value = -283.33334  # ALC!M2: -278.33334 -> -283.33334
$ws.Cells.Item(2, 14).Value = -458.5  # ALC!N2: None -> -458.5

$ws.Cells.Item(32, 8).Value = 1001  # ALC!H32: 0 -> 1001
$ws.Cells.Item(32, 9).Value = 1001  # ALC!I32: 0 -> 1001
$ws.Cells.Item(32, 11).Value = 1001  # ALC!K32: 0 -> 1001
$ws.Cells.Item(32, 13).Value = -675  # ALC!M32: None -> -675

$ws.Cells.Item(40, 8).Value = 1899.6  # ALC!H40: 2059.6 -> 1899.6
$ws.Cells.Item(40, 10).Value = 1299  # ALC!J40: 1899.3334 -> 1299
$ws.Cells.Item(40, 12).Value = 1299  # ALC!L40: 1899.3334 -> 1299
$ws.Cells.Item(40, 14).Value = -1649  # ALC!N40: -2249.3334 -> -1649

$ws.Cells.Item(88, 8).Value = 5001.5  # ALC!H88: 2041.8572 -> 5001.5
$ws.Cells.Item(88, 9).Value = 0  # ALC!I88: 1831.6666 -> 0
$ws.Cells.Item(88, 10).Value = 5001.5  # ALC!J88: 2199.5 -> 5001.5
$ws.Cells.Item(88, 11).Value = 0  # ALC!K88: 1831.6666 -> 0
$ws.Cells.Item(88, 12).Value = 5001.5  # ALC!L88: 2199.5 -> 5001.5
$ws.Cells.Item(88, 13).ClearContents()  # ALC!M88: -1425.6666 -> (deleted)
$ws.Cells.Item(88, 14).Value = -5813.5  # ALC!N88: -3011.5 -> -5813.5

$ws.Cells.Item(91, 8).Value = 5001.5  # ALC!H91: 2041.8572 -> 5001.5
$ws.Cells.Item(91, 9).Value = 0  # ALC!I91: 1831.6666 -> 0
$ws.Cells.Item(91, 10).Value = 5001.5  # ALC!J91: 2199.5 -> 5001.5
$ws.Cells.Item(91, 11).Value = 0  # ALC!K91: 1831.6666 -> 0
$ws.Cells.Item(91, 12).Value = 5001.5  # ALC!L91: 2199.5 -> 5001.5
$ws.Cells.Item(91, 13).ClearContents()  # ALC!M91: -427.6666 -> (deleted)
$ws.Cells.Item(91, 14).Value = -7809.5  # ALC!N91: -5007.5 -> -7809.5

$ws.Cells.Item(92, 8).Value = 919.6  # ALC!H92: 200000770 -> 919.6
$ws.Cells.Item(92, 9).Value = 799.3333  # ALC!I92: 333333900 -> 799.3333
$ws.Cells.Item(92, 11).Value = 799.3333  # ALC!K92: 333333900 -> 799.3333
$ws.Cells.Item(92, 13).Value = 448.6667  # ALC!M92: -333332652 -> 448.6667

$ws.Cells.Item(100, 8).Value = 750  # ALC!H100: 500 -> 750
$ws.Cells.Item(100, 10).Value = 1000  # ALC!J100: 0 -> 1000
$ws.Cells.Item(100, 12).Value = 1000  # ALC!L100: 0 -> 1000
$ws.Cells.Item(100, 14).Value = -2082  # ALC!N100: None -> -2082

$ws.Cells.Item(113, 8).Value = 13899.5  # ALC!H113: 0 -> 13899.5
$ws.Cells.Item(113, 9).Value = 13899.5  # ALC!I113: 0 -> 13899.5
$ws.Cells.Item(113, 11).Value = 13899.5  # ALC!K113: 0 -> 13899.5
$ws.Cells.Item(113, 13).Value = -10645.5  # ALC!M113: None -> -10645.5

$ws.Cells.Item(135, 8).Value = 1377.4286  # ALC!H135: 1114.4 -> 1377.4286
$ws.Cells.Item(135, 9).Value = 1128.4  # ALC!I135: 893 -> 1128.4
$ws.Cells.Item(135, 11).Value = 10155.6  # ALC!K135: 8037 -> 10155.6
$ws.Cells.Item(135, 13).Value = -7620.6  # ALC!M135: -5502 -> -7620.6

$ws.Cells.Item(138, 8).Value = 2549.4285  # ALC!H138: 2617.05 -> 2549.4285
$ws.Cells.Item(138, 9).Value = 1330.8334  # ALC!I138: 1357.6 -> 1330.8334
$ws.Cells.Item(138, 11).Value = 3992.5002  # ALC!K138: 4072.8 -> 3992.5002
$ws.Cells.Item(138, 13).Value = 1147.4998  # ALC!M138: 1067.2 -> 1147.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 302.27777  # ARM!H2: 302.94446 -> 302.27777
$ws.Cells.Item(2, 9).Value = 312.70587  # ARM!I2: 313.41177 -> 312.70587
$ws.Cells.Item(2, 11).Value = 312.70587  # ARM!K2: 313.41177 -> 312.70587
$ws.Cells.Item(2, 13).Value = -199.70587  # ARM!M2: -200.41177 -> -199.70587

$ws.Cells.Item(26, 8).Value = 760  # ARM!H26: 1700 -> 760
$ws.Cells.Item(26, 9).Value = 657.1429000000001  # ARM!I26: 725 -> 657.1429000000001
$ws.Cells.Item(26, 10).Value = 1000  # ARM!J26: 3000 -> 1000
$ws.Cells.Item(26, 11).Value = 657.1429000000001  # ARM!K26: 725 -> 657.1429000000001
$ws.Cells.Item(26, 12).Value = 1000  # ARM!L26: 3000 -> 1000
$ws.Cells.Item(26, 13).Value = -327.1429000000001  # ARM!M26: -395 -> -327.1429000000001
$ws.Cells.Item(26, 14).Value = -1660  # ARM!N26: -3660 -> -1660

$ws.Cells.Item(39, 8).Value = 6279.5713  # ARM!H39: 7223.5 -> 6279.5713
$ws.Cells.Item(39, 9).Value = 4811.4  # ARM!I39: 5860.25 -> 4811.4
$ws.Cells.Item(39, 11).Value = 4811.4  # ARM!K39: 5860.25 -> 4811.4
$ws.Cells.Item(39, 13).Value = -4291.4  # ARM!M39: -5340.25 -> -4291.4

$ws.Cells.Item(50, 8).Value = 15685.667  # ARM!H50: 15719 -> 15685.667
$ws.Cells.Item(50, 10).Value = 43855  # ARM!J50: 43955 -> 43855
$ws.Cells.Item(50, 12).Value = 43855  # ARM!L50: 43955 -> 43855
$ws.Cells.Item(50, 14).Value = -45283  # ARM!N50: -45383 -> -45283

$ws.Cells.Item(74, 8).Value = 2877.842  # ARM!H74: 2761.45 -> 2877.842
$ws.Cells.Item(74, 9).Value = 2936.7646  # ARM!I74: 2804.1667 -> 2936.7646
$ws.Cells.Item(74, 11).Value = 2936.7646  # ARM!K74: 2804.1667 -> 2936.7646
$ws.Cells.Item(74, 13).Value = -2062.7646  # ARM!M74: -1930.1667 -> -2062.7646

$ws.Cells.Item(77, 8).Value = 2877.842  # ARM!H77: 2761.45 -> 2877.842
$ws.Cells.Item(77, 9).Value = 2936.7646  # ARM!I77: 2804.1667 -> 2936.7646
$ws.Cells.Item(77, 11).Value = 14683.823  # ARM!K77: 14020.8335 -> 14683.823
$ws.Cells.Item(77, 13).Value = -10315.823  # ARM!M77: -9652.833500000001 -> -10315.823

$ws.Cells.Item(116, 8).Value = 302.27777  # ARM!H116: 302.94446 -> 302.27777
$ws.Cells.Item(116, 9).Value = 312.70587  # ARM!I116: 313.41177 -> 312.70587
$ws.Cells.Item(116, 11).Value = 312.70587  # ARM!K116: 313.41177 -> 312.70587
$ws.Cells.Item(116, 13).Value = 1981.29413  # ARM!M116: 1980.58823 -> 1981.29413

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 302.27777  # BSM!H3: 302.94446 -> 302.27777
$ws.Cells.Item(3, 9).Value = 312.70587  # BSM!I3: 313.41177 -> 312.70587
$ws.Cells.Item(3, 11).Value = 312.70587  # BSM!K3: 313.41177 -> 312.70587
$ws.Cells.Item(3, 13).Value = -198.70587  # BSM!M3: -199.41177 -> -198.70587

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1900  # CRP!H31: 550 -> 1900
$ws.Cells.Item(31, 9).Value = 0  # CRP!I31: 550 -> 0
$ws.Cells.Item(31, 10).Value = 1900  # CRP!J31: 0 -> 1900
$ws.Cells.Item(31, 11).Value = 0  # CRP!K31: 550 -> 0
$ws.Cells.Item(31, 12).Value = 1900  # CRP!L31: 0 -> 1900
$ws.Cells.Item(31, 13).ClearContents()  # CRP!M31: -255 -> (deleted)
$ws.Cells.Item(31, 14).Value = -2490  # CRP!N31: None -> -2490

$ws.Cells.Item(34, 8).Value = 1900  # CRP!H34: 550 -> 1900
$ws.Cells.Item(34, 9).Value = 0  # CRP!I34: 550 -> 0
$ws.Cells.Item(34, 10).Value = 1900  # CRP!J34: 0 -> 1900
$ws.Cells.Item(34, 11).Value = 0  # CRP!K34: 550 -> 0
$ws.Cells.Item(34, 12).Value = 1900  # CRP!L34: 0 -> 1900
$ws.Cells.Item(34, 13).ClearContents()  # CRP!M34: -348 -> (deleted)
$ws.Cells.Item(34, 14).Value = -2304  # CRP!N34: None -> -2304

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 2547.0293  # CUL!H4: 2444.0293 -> 2547.0293
$ws.Cells.Item(4, 9).Value = 1933.2858  # CUL!I4: 1868.0454 -> 1933.2858
$ws.Cells.Item(4, 10).Value = 3538.4614  # CUL!J4: 3500 -> 3538.4614
$ws.Cells.Item(4, 11).Value = 5799.857400000001  # CUL!K4: 5604.1362 -> 5799.857400000001
$ws.Cells.Item(4, 12).Value = 10615.3842  # CUL!L4: 10500 -> 10615.3842
$ws.Cells.Item(4, 13).Value = -5687.857400000001  # CUL!M4: -5492.1362 -> -5687.857400000001
$ws.Cells.Item(4, 14).Value = -10839.3842  # CUL!N4: -10724 -> -10839.3842

$ws.Cells.Item(23, 8).Value = 118.5  # CUL!H23: 112.2 -> 118.5
$ws.Cells.Item(23, 10).Value = 102.2  # CUL!J23: 90.25 -> 102.2
$ws.Cells.Item(23, 12).Value = 306.6  # CUL!L23: 270.75 -> 306.6
$ws.Cells.Item(23, 14).Value = -776.6  # CUL!N23: -740.75 -> -776.6

$ws.Cells.Item(31, 8).Value = 495  # CUL!H31: 0 -> 495
$ws.Cells.Item(31, 10).Value = 495  # CUL!J31: 0 -> 495
$ws.Cells.Item(31, 12).Value = 1485  # CUL!L31: 0 -> 1485
$ws.Cells.Item(31, 14).Value = -2061  # CUL!N31: None -> -2061

$ws.Cells.Item(50, 8).Value = 1264.9166  # CUL!H50: 1348.2727 -> 1264.9166
$ws.Cells.Item(50, 9).Value = 382.7143  # CUL!I50: 388.5 -> 382.7143
$ws.Cells.Item(50, 11).Value = 1148.1429  # CUL!K50: 1165.5 -> 1148.1429
$ws.Cells.Item(50, 13).Value = -667.1428999999998  # CUL!M50: -684.5 -> -667.1428999999998

$ws.Cells.Item(53, 8).Value = 1264.9166  # CUL!H53: 1348.2727 -> 1264.9166
$ws.Cells.Item(53, 9).Value = 382.7143  # CUL!I53: 388.5 -> 382.7143
$ws.Cells.Item(53, 11).Value = 1148.1429  # CUL!K53: 1165.5 -> 1148.1429
$ws.Cells.Item(53, 13).Value = -667.1428999999998  # CUL!M53: -684.5 -> -667.1428999999998

$ws.Cells.Item(55, 8).Value = 2007  # CUL!H55: 2158.3333 -> 2007
$ws.Cells.Item(55, 9).Value = 1429.9  # CUL!I55: 1466.6666 -> 1429.9
$ws.Cells.Item(55, 10).Value = 3449.75  # CUL!J55: 4233.3335 -> 3449.75
$ws.Cells.Item(55, 11).Value = 4289.700000000001  # CUL!K55: 4399.9998 -> 4289.700000000001
$ws.Cells.Item(55, 12).Value = 10349.25  # CUL!L55: 12700.0005 -> 10349.25
$ws.Cells.Item(55, 13).Value = -4112.700000000001  # CUL!M55: -4222.9998 -> -4112.700000000001
$ws.Cells.Item(55, 14).Value = -10703.25  # CUL!N55: -13054.0005 -> -10703.25

$ws.Cells.Item(69, 8).Value = 6000  # CUL!H69: 0 -> 6000
$ws.Cells.Item(69, 10).Value = 6000  # CUL!J69: 0 -> 6000
$ws.Cells.Item(69, 12).Value = 18000  # CUL!L69: 0 -> 18000
$ws.Cells.Item(69, 14).Value = -19622  # CUL!N69: None -> -19622

$ws.Cells.Item(72, 8).Value = 6000  # CUL!H72: 0 -> 6000
$ws.Cells.Item(72, 10).Value = 6000  # CUL!J72: 0 -> 6000
$ws.Cells.Item(72, 12).Value = 54000  # CUL!L72: 0 -> 54000
$ws.Cells.Item(72, 14).Value = -62112  # CUL!N72: None -> -62112

$ws.Cells.Item(93, 8).Value = 1800  # CUL!H93: 1766.6666 -> 1800
$ws.Cells.Item(93, 10).Value = 1800  # CUL!J93: 1766.6666 -> 1800
$ws.Cells.Item(93, 12).Value = 5400  # CUL!L93: 5299.9998 -> 5400
$ws.Cells.Item(93, 14).Value = -9144  # CUL!N93: -9043.9998 -> -9144

$ws.Cells.Item(121, 8).Value = 1033  # CUL!H121: 370 -> 1033
$ws.Cells.Item(121, 9).Value = 0  # CUL!I121: 370 -> 0
$ws.Cells.Item(121, 10).Value = 1033  # CUL!J121: 0 -> 1033
$ws.Cells.Item(121, 11).Value = 0  # CUL!K121: 1110 -> 0
$ws.Cells.Item(121, 12).Value = 3099  # CUL!L121: 0 -> 3099
$ws.Cells.Item(121, 13).ClearContents()  # CUL!M121: 200 -> (deleted)
$ws.Cells.Item(121, 14).Value = -5719  # CUL!N121: None -> -5719

$ws.Cells.Item(129, 8).Value = 997.55554  # CUL!H129: 1146.091 -> 997.55554
$ws.Cells.Item(129, 9).Value = 630  # CUL!I129: 496.33334 -> 630
$ws.Cells.Item(129, 10).Value = 1102.5714  # CUL!J129: 1389.75 -> 1102.5714
$ws.Cells.Item(129, 11).Value = 1890  # CUL!K129: 1489.00002 -> 1890
$ws.Cells.Item(129, 12).Value = 3307.7142  # CUL!L129: 4169.25 -> 3307.7142
$ws.Cells.Item(129, 13).Value = 3110  # CUL!M129: 3510.99998 -> 3110
$ws.Cells.Item(129, 14).Value = -13307.7142  # CUL!N129: -14169.25 -> -13307.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1127  # GSM!H132: 1095 -> 1127
$ws.Cells.Item(132, 9).Value = 1127  # GSM!I132: 1095 -> 1127
$ws.Cells.Item(132, 11).Value = 3381  # GSM!K132: 3285 -> 3381
$ws.Cells.Item(132, 13).Value = -851  # GSM!M132: -755 -> -851

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1757.5  # LTW!H16: 2212.5 -> 1757.5
$ws.Cells.Item(16, 9).Value = 1619.4445  # LTW!I16: 2028.5714 -> 1619.4445
$ws.Cells.Item(16, 10).Value = 3000  # LTW!J16: 3500 -> 3000
$ws.Cells.Item(16, 11).Value = 1619.4445  # LTW!K16: 2028.5714 -> 1619.4445
$ws.Cells.Item(16, 12).Value = 3000  # LTW!L16: 3500 -> 3000
$ws.Cells.Item(16, 13).Value = -1449.4445  # LTW!M16: -1858.5714 -> -1449.4445
$ws.Cells.Item(16, 14).Value = -3340  # LTW!N16: -3840 -> -3340

$ws.Cells.Item(22, 8).Value = 4999.4287  # LTW!H22: 5000.5 -> 4999.4287
$ws.Cells.Item(22, 10).Value = 4999  # LTW!J22: 0 -> 4999
$ws.Cells.Item(22, 12).Value = 4999  # LTW!L22: 0 -> 4999
$ws.Cells.Item(22, 14).Value = -5589  # LTW!N22: None -> -5589

$ws.Cells.Item(27, 8).Value = 4999.4287  # LTW!H27: 5000.5 -> 4999.4287
$ws.Cells.Item(27, 10).Value = 4999  # LTW!J27: 0 -> 4999
$ws.Cells.Item(27, 12).Value = 4999  # LTW!L27: 0 -> 4999
$ws.Cells.Item(27, 14).Value = -5213  # LTW!N27: None -> -5213

$ws.Cells.Item(43, 8).Value = 0  # LTW!H43: 11000 -> 0
$ws.Cells.Item(43, 10).Value = 0  # LTW!J43: 11000 -> 0
$ws.Cells.Item(43, 12).Value = 0  # LTW!L43: 11000 -> 0
$ws.Cells.Item(43, 14).ClearContents()  # LTW!N43: -11386 -> (deleted)

$ws.Cells.Item(46, 8).Value = 2000000  # LTW!H46: 503243.75 -> 2000000
$ws.Cells.Item(46, 9).Value = 2000000  # LTW!I46: 1000750 -> 2000000
$ws.Cells.Item(46, 10).Value = 0  # LTW!J46: 5737.5 -> 0
$ws.Cells.Item(46, 11).Value = 2000000  # LTW!K46: 1000750 -> 2000000
$ws.Cells.Item(46, 12).Value = 0  # LTW!L46: 5737.5 -> 0
$ws.Cells.Item(46, 13).Value = -1999812  # LTW!M46: -1000562 -> -1999812
$ws.Cells.Item(46, 14).ClearContents()  # LTW!N46: -6113.5 -> (deleted)

$ws.Cells.Item(48, 8).Value = 0  # LTW!H48: 27845 -> 0
$ws.Cells.Item(48, 10).Value = 0  # LTW!J48: 27845 -> 0
$ws.Cells.Item(48, 12).Value = 0  # LTW!L48: 27845 -> 0
$ws.Cells.Item(48, 14).ClearContents()  # LTW!N48: -29167 -> (deleted)

$ws.Cells.Item(55, 8).Value = 608.3333  # LTW!H55: 660.7143 -> 608.3333
$ws.Cells.Item(55, 9).Value = 520  # LTW!I55: 583.3333 -> 520
$ws.Cells.Item(55, 11).Value = 520  # LTW!K55: 583.3333 -> 520
$ws.Cells.Item(55, 13).Value = -347  # LTW!M55: -410.3333 -> -347

$ws.Cells.Item(82, 8).Value = 1529.1428  # LTW!H82: 1587 -> 1529.1428
$ws.Cells.Item(82, 9).Value = 1390.8  # LTW!I82: 1443 -> 1390.8
$ws.Cells.Item(82, 11).Value = 1390.8  # LTW!K82: 1443 -> 1390.8
$ws.Cells.Item(82, 13).Value = -1029.8  # LTW!M82: -1082 -> -1029.8

$ws.Cells.Item(85, 8).Value = 1529.1428  # LTW!H85: 1587 -> 1529.1428
$ws.Cells.Item(85, 9).Value = 1390.8  # LTW!I85: 1443 -> 1390.8
$ws.Cells.Item(85, 11).Value = 1390.8  # LTW!K85: 1443 -> 1390.8
$ws.Cells.Item(85, 13).Value = -142.8  # LTW!M85: -195 -> -142.8

$ws.Cells.Item(132, 8).Value = 0  # LTW!H132: 2099.5 -> 0
$ws.Cells.Item(132, 9).Value = 0  # LTW!I132: 2099.5 -> 0
$ws.Cells.Item(132, 11).Value = 0  # LTW!K132: 6298.5 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # LTW!M132: -3768.5 -> (deleted)

$ws.Cells.Item(136, 8).Value = 0  # LTW!H136: 996 -> 0
$ws.Cells.Item(136, 9).Value = 0  # LTW!I136: 996 -> 0
$ws.Cells.Item(136, 11).Value = 0  # LTW!K136: 2988 -> 0
$ws.Cells.Item(136, 13).ClearContents()  # LTW!M136: -438 -> (deleted)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 5500000  # WVR!H5: 4503333.5 -> 5500000
$ws.Cells.Item(5, 9).Value = 10000000  # WVR!I5: 5005000 -> 10000000
$ws.Cells.Item(5, 10).Value = 4000000  # WVR!J5: 4252500 -> 4000000
$ws.Cells.Item(5, 11).Value = 10000000  # WVR!K5: 5005000 -> 10000000
$ws.Cells.Item(5, 12).Value = 4000000  # WVR!L5: 4252500 -> 4000000
$ws.Cells.Item(5, 13).Value = -9999888  # WVR!M5: -5004888 -> -9999888
$ws.Cells.Item(5, 14).Value = -4000224  # WVR!N5: -4252724 -> -4000224

$ws.Cells.Item(136, 8).Value = 973  # WVR!H136: 1550.1 -> 973
$ws.Cells.Item(136, 9).Value = 973  # WVR!I136: 1389 -> 973
$ws.Cells.Item(136, 10).Value = 0  # WVR!J136: 3000 -> 0
$ws.Cells.Item(136, 11).Value = 2919  # WVR!K136: 4167 -> 2919
$ws.Cells.Item(136, 12).Value = 0  # WVR!L136: 9000 -> 0
$ws.Cells.Item(136, 13).Value = -369  # WVR!M136: -1617 -> -369
$ws.Cells.Item(136, 14).ClearContents()  # WVR!N136: -14100 -> (deleted)

$ws.Cells.Item(140, 8).Value = 60000  # WVR!H140: 65000 -> 60000
$ws.Cells.Item(140, 10).Value = 60000  # WVR!J140: 65000 -> 60000
$ws.Cells.Item(140, 12).Value = 60000  # WVR!L140: 65000 -> 60000
$ws.Cells.Item(140, 14).Value = -70360  # WVR!N140: -75360 -> -70360
